# Append the newest bitcoin buy entry after running on 2026-01-04
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Column A: date stored as plain text (matches the sheet's existing convention
# of recording later dates as literal text strings rather than date serials).
# Force text entry via the "@" number format, then reset the cell style back
# to Normal so no extra style index is left behind on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01/04/2026"
$ws.Cells.Item($row, 1).Style = "Normal"

# Column B: coins purchased
$ws.Cells.Item($row, 2).Value = 0.0005391099999999989

# Column C: price per coin
$ws.Cells.Item($row, 3).Value = 91817.99632728033

# Column D: cost
$ws.Cells.Item($row, 4).Value = 50
